$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'49.695.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").Value = "'2.615.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.29%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'325.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "'110.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").Value = "'0.998"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "'0.557"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("D10").Value = "'40.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("D11").Value = "'20.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").Value = "'0.0818"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D15").Value = "'3.021.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.23%  "
$ws.Range("D16").Value = "'2.613.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.18%  "
$ws.Range("D17").Value = "'0.871"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("D18").Value = "'49.560.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("D19").Value = "'3.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.75%  "
$ws.Range("D20").Value = "'13.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("D21").Value = "'6.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'0.0₃0952"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").Value = "'72.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "'278.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "'26.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "'2.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("D29").Value = "'9.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.74%  "
$ws.Range("E30").Value = "  +2.24%  "
$ws.Range("D31").Value = "'36.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.83%  "
$ws.Range("D32").Value = "'49.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "'19.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "'5.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("D36").Value = "'0.0791"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("E37").Value = "  +4.62%  "
$ws.Range("D38").Value = "'4.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("E39").Value = "  +5.80%  "
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").Value = "'122.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("D42").Value = "'22.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.82%  "
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").Value = "'0.0315"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.95%  "
$ws.Range("D45").Value = "'3.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.54%  "
$ws.Range("D46").Value = "'2.046.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("D47").Value = "'2.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.55%  "
$ws.Range("E48").Value = "  +8.72%  "
$ws.Range("D49").Value = "'9.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  +3.06%  "
$ws.Range("D51").Value = "'81.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.07%  "
